$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-formatted numbers in the source data (e.g. "557.99" as a
# literal string, not a number). Where the new value would otherwise be auto-parsed by
# Excel as a numeric value (losing formatting like trailing zeros), force the cell to
# Text format first so the literal string is preserved, matching the source workbook.

$ws.Range("D2").Value = "62.998.07"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "3.024.83"
$ws.Range("E3").Value = "  -3.46%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.99"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.48"
$ws.Range("E6").Value = "  -5.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").Value = "  -2.48%  "

$ws.Range("D9").Value = "3.026.35"
$ws.Range("E9").Value = "  -3.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.113"
$ws.Range("E10").Value = "  -1.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.35"
$ws.Range("E11").Value = "  -4.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.367"
$ws.Range("E12").Value = "  -2.87%  "

$ws.Range("D13").Value = "3.548.59"
$ws.Range("E13").Value = "  -3.41%  "

$ws.Range("E14").Value = "  -3.07%  "

$ws.Range("D15").Value = "63.016.79"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.93"
$ws.Range("E16").Value = "  -3.58%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000150"
$ws.Range("E17").Value = "  -1.82%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.011.59"
$ws.Range("E18").Value = "  -3.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "397.26"
$ws.Range("E19").Value = "  -1.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.13"
$ws.Range("E20").Value = "  -1.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.97"
$ws.Range("E21").Value = "  -3.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.67"
$ws.Range("E22").Value = "  -4.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.25"
$ws.Range("E24").Value = "  -2.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.468"
$ws.Range("E25").Value = "  -2.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.187"
$ws.Range("E26").Value = "  -6.61%  "

$ws.Range("D27").Value = "0.0₃0977"
$ws.Range("E27").Value = "  -2.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.63"
$ws.Range("E28").Value = "  -0.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.75"
$ws.Range("E31").Value = "  -1.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.59"
$ws.Range("E32").Value = "  -1.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "161.06"
$ws.Range("E33").Value = "  +5.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.72"
$ws.Range("E34").Value = "  -0.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.05"
$ws.Range("E35").Value = "  -2.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.10"
$ws.Range("E36").Value = "  -0.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.30"
$ws.Range("E37").Value = "  -1.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.59"
$ws.Range("E38").Value = "  -3.42%  "

$ws.Range("D39").Value = "2.481.20"
$ws.Range("E39").Value = "  -9.90%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.69"
$ws.Range("E40").Value = "  -2.64%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.94"
$ws.Range("E41").Value = "  -2.19%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.75"
$ws.Range("E42").Value = "  -2.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.667"
$ws.Range("E43").Value = "  -3.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0598"
$ws.Range("E44").Value = "  -2.99%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0249"
$ws.Range("E45").Value = "  -2.71%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.996"
$ws.Range("E46").Value = "  -0.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.99"
$ws.Range("E47").Value = "  -6.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.96"
$ws.Range("E48").Value = "  -3.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0952"
$ws.Range("E49").Value = "  -1.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.47"
$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "264.48"
$ws.Range("E51").Value = "  -6.06%  "
